# Add a new "2023" column (L) to the table, mirroring the existing 2022
# column (K): a blank bordered cell in the separator row, a right-aligned
# year header, and the new data value. Row 5 also grows a touch to fit the
# extra wrapped header text, and the stale saved cell selection is reset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 is a blank bottom-bordered separator row spanning the table; extend
# it to column L by cloning K3's formatting (style index 4).
$ws.Range("K3").Copy()
$ws.Range("L3").PasteSpecial(-4122)

# Header cell L4 = 2023, cloning K4's formatting (style index 9).
$ws.Range("K4").Copy()
$ws.Range("L4").PasteSpecial(-4122)
$ws.Range("L4").Value = 2023

# Data cell L5 = R&D expenditure share of GDP for 2023, cloning K5's
# formatting (style index 13).
$ws.Range("K5").Copy()
$ws.Range("L5").PasteSpecial(-4122)
$ws.Range("L5").Value = 0.11972285283622097

# Row 5 grows from 36.75 to 40.5 points to accommodate the new column.
$ws.Rows.Item(5).RowHeight = 40.5

# Reset the lingering cell selection left over from the previous edit
# session back to the top-left cell.
$ws.Range("A1").Select()
